# Update "想去人数" (interest count) figures to match the newer scrape
# snapshot (gh-pages output regenerated at commit 456a3b4).
#
# The workbook has 4 sheets:
#   展览     (Exhibitions)      - sheet 1
#   演出     (Performances)     - sheet 2
#   本地生活 (Local life)       - sheet 3 (unchanged - header row only)
#   全部类型 (All types, union) - sheet 4 (mirrors sheet1+sheet2 rows)
#
# Column F on each data sheet holds the interest-count value that was
# refreshed by the scraper; only those cells changed between commits.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 20623
$ws1.Range("F4").Value  = 327
$ws1.Range("F5").Value  = 1111
$ws1.Range("F6").Value  = 21
$ws1.Range("F7").Value  = 7721
$ws1.Range("F8").Value  = 536
$ws1.Range("F9").Value  = 748
$ws1.Range("F11").Value = 51
$ws1.Range("F16").Value = 209
$ws1.Range("F18").Value = 479
$ws1.Range("F20").Value = 694
$ws1.Range("F24").Value = 337
$ws1.Range("F25").Value = 1151
$ws1.Range("F30").Value = 582
$ws1.Range("F32").Value = 4923
$ws1.Range("F35").Value = 55
$ws1.Range("F36").Value = 12846
$ws1.Range("F37").Value = 1348
$ws1.Range("F38").Value = 103
$ws1.Range("F39").Value = 39
$ws1.Range("F43").Value = 4027

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 237

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 20623
$ws4.Range("F4").Value  = 327
$ws4.Range("F5").Value  = 1111
$ws4.Range("F6").Value  = 21
$ws4.Range("F7").Value  = 7721
$ws4.Range("F8").Value  = 536
$ws4.Range("F9").Value  = 748
$ws4.Range("F11").Value = 51
$ws4.Range("F16").Value = 209
$ws4.Range("F18").Value = 479
$ws4.Range("F20").Value = 694
$ws4.Range("F24").Value = 337
$ws4.Range("F25").Value = 1151
$ws4.Range("F29").Value = 237
$ws4.Range("F31").Value = 582
$ws4.Range("F35").Value = 4923
$ws4.Range("F38").Value = 55
$ws4.Range("F39").Value = 12846
$ws4.Range("F40").Value = 1348
$ws4.Range("F41").Value = 103
$ws4.Range("F42").Value = 39
$ws4.Range("F46").Value = 4027
